# Generate Report for handoff
# Update "Latest Handoff Datetime" (column D) for the rows that share the
# placeholder handoff timestamp, on both the zh-cn and de-de sheets, to
# reflect the newly generated handoff.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcnRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $zhcnRows) {
    $zhcn.Cells.Item($r, 4).Value = "2016-02-17 03:42:10"
}

$dedeRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $dedeRows) {
    $dede.Cells.Item($r, 4).Value = "2016-02-17 03:42:21"
}
